$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.48109986927249
$ws.Range("B3").Value = 3.81979832202712
$ws.Range("B4").Value = 2.05995993267184
$ws.Range("B5").Value = 0.795682877021712
$ws.Range("B6").Value = 2.22386194507141
$ws.Range("B7").Value = 2.07608058921291
